$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"   = -8.183
    "D14"  = -7.804
    "D21"  = -8.227
    "D23"  = -7.980000000000001
    "D25"  = -7.842999999999999
    "D26"  = -7.399000000000001
    "D29"  = -7.137
    "D53"  = -7.947999999999999
    "D57"  = -8
    "D59"  = -8.234
    "D69"  = -7.640000000000001
    "D79"  = -7.98
    "D83"  = -8.125
    "D91"  = -7.441
    "D93"  = -6.976999999999999
    "D103" = -8.211000000000002
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
